$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 1.925925925925943
$ws.Range("N2").Value = 1.459904774678112
$ws.Range("O2").Value = 1.554373915558126

$ws.Range("I5").Value = 1.925925925925943
$ws.Range("N5").Value = 1.459904774678112
$ws.Range("O5").Value = 1.554373915558126

$ws.Range("I9").Value = 13.17361111111111
$ws.Range("N9").Value = 1.58937742977605
$ws.Range("O9").Value = 1.704024252511443

$ws.Range("I10").Value = 13.17361111111111
$ws.Range("N10").Value = 1.58937742977605
$ws.Range("O10").Value = 1.704024252511443
